# Update model3_df_results.xlsx: refresh R^2 (C), RMSE (D) and U (E) metrics
# for rows 2-9, and recolor the RMSE/U heat-map cells to match the new
# background-gradient coloring that corresponds to the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FillColor {
    param($range, [int]$r, [int]$g, [int]$b)
    $range.Interior.Color = $r + ($g * 256) + ($b * 65536)
}

# Row 2
$ws.Range("C2").Value = -30.0126
$ws.Range("D2").Value = 1.5189
Set-FillColor $ws.Range("D2") 247 252 245
$ws.Range("E2").Value = 4.685
Set-FillColor $ws.Range("E2") 247 252 245

# Row 3
$ws.Range("C3").Value = -11.4739
$ws.Range("D3").Value = 1.4664
Set-FillColor $ws.Range("D3") 186 227 179
$ws.Range("E3").Value = 3.4766
Set-FillColor $ws.Range("E3") 12 119 53

# Row 4
$ws.Range("C4").Value = -5.6612
$ws.Range("D4").Value = 1.4309
Set-FillColor $ws.Range("D4") 117 196 119
$ws.Range("E4").Value = 3.437
Set-FillColor $ws.Range("E4") 5 113 47

# Row 5
$ws.Range("C5").Value = -3.0416
$ws.Range("D5").Value = 1.3929
Set-FillColor $ws.Range("D5") 44 148 76
$ws.Range("E5").Value = 3.3217
Set-FillColor $ws.Range("E5") 0 89 36

# Row 6
$ws.Range("C6").Value = -2.2067
$ws.Range("D6").Value = 1.3874
Set-FillColor $ws.Range("D6") 36 140 70
$ws.Range("E6").Value = 3.3197
Set-FillColor $ws.Range("E6") 0 89 36

# Row 7
$ws.Range("C7").Value = -1.9761
$ws.Range("D7").Value = 1.3419
Set-FillColor $ws.Range("D7") 0 68 27
$ws.Range("E7").Value = 3.2641
Set-FillColor $ws.Range("E7") 0 76 30

# Row 8
$ws.Range("C8").Value = -1.8975
$ws.Range("D8").Value = 1.3462
Set-FillColor $ws.Range("D8") 0 76 30
$ws.Range("E8").Value = 3.2768
Set-FillColor $ws.Range("E8") 0 78 31

# Row 9
$ws.Range("C9").Value = -1.8189
$ws.Range("D9").Value = 1.3541
Set-FillColor $ws.Range("D9") 0 90 36
$ws.Range("E9").Value = 3.2277
Set-FillColor $ws.Range("E9") 0 68 27
